$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2: report-type code changes from 001 to 004 (keep as text, preserve leading zeros)
$ws.Range("J2").Value = "'004"

# N2: report date moves forward a quarter
$ws.Range("N2").Value = "2020-09-30 00:00:00"

# Updated financial figures for the new report date
$ws.Range("O2").Value = 705520703.62
$ws.Range("P2").Value = 118948571.13
$ws.Range("Q2").Value = 72927195.19
$ws.Range("S2").Value = 250763856.87
$ws.Range("U2").Value = 135221919.69
$ws.Range("W2").Value = 300943890.59
$ws.Range("X2").Value = 76633129.23999999
$ws.Range("AB2").Value = 404576813.03
$ws.Range("AF2").Value = 188.5583848424
$ws.Range("AG2").Value = 42.6555718416

# Ratio columns that no longer have data for this report - clear them out
$ws.Range("R2").ClearContents()
$ws.Range("T2").ClearContents()
$ws.Range("V2").ClearContents()
$ws.Range("Y2").ClearContents()
$ws.Range("Z2").ClearContents()
$ws.Range("AA2").ClearContents()
$ws.Range("AC2").ClearContents()
$ws.Range("AD2").ClearContents()
$ws.Range("AE2").ClearContents()
